$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.338.49'
$ws.Range("E2").Value = '  -0.65%  '

$ws.Range("D3").Value = '1.869.53'
$ws.Range("E3").Value = '  -0.23%  '

$ws.Range("E4").Value = '  -0.04%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '243.60'
$ws.Range("E5").Value = '  -1.58%  '

$ws.Range("E6").Value = '  -0.03%  '

$ws.Range("E7").Value = '  -1.16%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2870'
$ws.Range("E8").Value = '  -1.35%  '

$ws.Range("E9").Value = '  -0.81%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '22.03'
$ws.Range("E10").Value = '  +0.85%  '

$ws.Range("E11").Value = '  +0.26%  '

$ws.Range("D12").Value = '1.870.25'
$ws.Range("E12").Value = '  -0.19%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '95.94'
$ws.Range("E13").Value = '  -0.06%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.7200'
$ws.Range("E14").Value = '  -2.28%  '

$ws.Range("E15").Value = '  -1.18%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '278.21'
$ws.Range("E16").Value = '  +1.52%  '

$ws.Range("D17").Value = '30.321.61'
$ws.Range("E17").Value = '  -0.89%  '

$ws.Range("E18").Value = '  -1.75%  '

$ws.Range("B19").Value = 'Dai'
$ws.Range("C19").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.000'
$ws.Range("E19").Value = '  -0.04%  '

$ws.Range("B20").Value = 'ShibaInu'
$ws.Range("C20").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.000007498'
$ws.Range("E20").Value = '  -0.03%  '

$ws.Range("D21").Value = '2.114.02'
$ws.Range("E21").Value = '  -0.26%  '

$ws.Range("E22").Value = '  -0.02%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.214'
$ws.Range("E23").Value = '  -0.03%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.221'
$ws.Range("E24").Value = '  +0.87%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '163.29'
$ws.Range("E25").Value = '  -0.94%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.029'
$ws.Range("E26").Value = '  -1.62%  '

$ws.Range("E27").Value = '  -0.72%  '

$ws.Range("E28").Value = '  -2.01%  '

$ws.Range("E29").Value = '  -1.32%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.474'
$ws.Range("E31").Value = '  -1.53%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.198'
$ws.Range("E32").Value = '  -1.28%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.088'
$ws.Range("E33").Value = '  +0.09%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.04799'
$ws.Range("E34").Value = '  +0.25%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.117'
$ws.Range("E35").Value = '  -0.30%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.6880'
$ws.Range("E36").Value = '  -0.79%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.711'
$ws.Range("E37").Value = '  -0.24%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01873'
$ws.Range("E38").Value = '  +1.14%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.807'
$ws.Range("E39").Value = '  +1.90%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.199'
$ws.Range("E40").Value = '  -1.31%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '74.13'
$ws.Range("E41").Value = '  +1.07%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.930'
$ws.Range("E42").Value = '  -2.75%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.4207'
$ws.Range("E43").Value = '  +0.37%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.9993'

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.8281'
$ws.Range("E45").Value = '  -0.83%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '100.48'
$ws.Range("E46").Value = '  -1.06%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '9.543'
$ws.Range("E47").Value = '  +2.11%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '35.21'
$ws.Range("E48").Value = '  -0.12%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '6.940'
$ws.Range("E49").Value = '  -0.25%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '899.86'
$ws.Range("E50").Value = '  -1.49%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.05717'
$ws.Range("E51").Value = '  +0.87%  '
